$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The "0" column (AK) was a leftover/duplicate column that held a constant
# value of 1 for every data row. Removing it shifts the "predicted" column
# (formerly AL) one position to the left, into AK.
$ws.Range("AK1").EntireColumn.Delete()
